$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.758.39"
$ws.Range("E2").Value = "  +0.27%  "
$ws.Range("D3").Value = "1.753.61"
$ws.Range("E3").Value = "  +1.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.07"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5060"
$ws.Range("E7").Value = "  +3.31%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "40.59"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.2625"
$ws.Range("E9").Value = "  +8.61%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06204"
$ws.Range("E10").Value = "  +3.31%  "
$ws.Range("D11").Value = "1.758.58"
$ws.Range("E11").Value = "  +1.30%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.06943"
$ws.Range("E12").Value = "  +3.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.47"
$ws.Range("E13").Value = "  +4.57%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6039"
$ws.Range("E14").Value = "  +2.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "78.26"
$ws.Range("E15").Value = "  +2.33%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.455"
$ws.Range("E16").Value = "  +1.75%  "
$ws.Range("E17").Value = "  +0.08%  "
$ws.Range("E18").Value = "  +0.00%  "
$ws.Range("D19").Value = "25.803.22"
$ws.Range("E19").Value = "  +0.26%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.66"
$ws.Range("E20").Value = "  +2.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006767"
$ws.Range("E21").Value = "  +6.96%  "
$ws.Range("D22").Value = "1.974.49"
$ws.Range("E22").Value = "  +0.81%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.057"
$ws.Range("E23").Value = "  +3.81%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.177"
$ws.Range("E24").Value = "  +4.02%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.178"
$ws.Range("E25").Value = "  +1.29%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "137.80"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.469"
$ws.Range("E27").Value = "  +2.89%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.09"
$ws.Range("E28").Value = "  +5.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.800"
$ws.Range("E29").Value = "  -2.35%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.45"
$ws.Range("E30").Value = "  +1.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08280"
$ws.Range("E31").Value = "  +1.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.698"
$ws.Range("E32").Value = "  +1.70%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.394"
$ws.Range("E33").Value = "  +2.75%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04371"
$ws.Range("E34").Value = "  +0.30%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.000"
$ws.Range("E35").Value = "  +0.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.646"
$ws.Range("E36").Value = "  -1.13%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.002"
$ws.Range("E37").Value = "  -2.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.6004"
$ws.Range("E38").Value = "  -1.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.697"
$ws.Range("E39").Value = "  -2.92%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.958"
$ws.Range("E40").Value = "  -5.19%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01547"
$ws.Range("E41").Value = "  +4.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.001"
$ws.Range("E42").Value = "  +0.19%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "103.39"
$ws.Range("E43").Value = "  +1.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7482"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.3803"
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.882"
$ws.Range("E46").Value = "  -4.58%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05486"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.1078"
$ws.Range("E48").Value = "  +4.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "5.949"
$ws.Range("E49").Value = "  -1.80%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.18"
$ws.Range("E50").Value = "  +1.61%  "
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.002"
$ws.Range("E51").Value = "  -0.05%  "
